$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, [string]$text)
    # Force the cell to stay a text/string cell even when the text looks
    # like a pure number (e.g. "1545"), mirroring how the source workbook
    # stores booking/room numbers as shared strings rather than numbers.
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

# --- Update existing status column (D) values ---
# Row 1 (header row's status sample) and row 2: CHECK-IN -> CHECK-OUT
$ws.Cells.Item(1,4).Value = "CHECK-OUT"
$ws.Cells.Item(2,4).Value = "CHECK-OUT"
# Row 3 stays CHECK-OUT (unchanged)
# Row 4: Waiting -> CHECK-IN
$ws.Cells.Item(4,4).Value = "CHECK-IN"

# --- Append two new booking rows (5 and 6) ---
Set-TextValue $ws.Cells.Item(5,1) "1545"
$ws.Cells.Item(5,2).Value = "Raweeroj   Thongdee"
Set-TextValue $ws.Cells.Item(5,3) "1003"
$ws.Cells.Item(5,4).Value = "Waiting"
$ws.Cells.Item(5,5).Value = "24-04-2020 03:32:19"

Set-TextValue $ws.Cells.Item(6,1) "1573"
$ws.Cells.Item(6,2).Value = "Raweeroj   Thongdee"
Set-TextValue $ws.Cells.Item(6,3) "2003"
$ws.Cells.Item(6,4).Value = "Waiting"
$ws.Cells.Item(6,5).Value = "24-04-2020 03:32:19"

# --- Widen status column (D) to fit the new "CHECK-OUT" content ---
$ws.Columns.Item(4).ColumnWidth = 10.833333333333334
